$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" date column (C) for rows 2-18 from 2023-10-05 (45204)
# to 2023-10-08 (45207), as reflected in the source data refresh.
for ($row = 2; $row -le 18; $row++) {
    $ws.Cells.Item($row, 3).Value = 45207
}
